$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the English (column C) translations: "version" -> "tag"
$ws.Range("C3").Value = "No images of a specified tag"
$ws.Range("C7").Value = "Delete images of a specified tag"
$ws.Range("C9").Value = "Delete images according to tag marker status"
$ws.Range("C10").Value = "Delete all images without tag markers"
$ws.Range("C11").Value = "Delete all images with tag markers"

# Apply a new (red Times New Roman, wrapped) style to column D
$colD = $ws.Columns.Item(4)
$colD.Font.Name = "Times New Roman"
$colD.Font.Size = 12
$colD.Font.Color = 255
$colD.WrapText = $true

# Move the selection / active cell like the author's session ended up
$ws.Range("C17").Select()
